# Update workbook to add data for 2022-09-23 (new "through" date: September 15, 2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the (only) worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-09-15"

# Update the column-B header text (shared string) to reflect the new "through" date
$ws.Range("B1").Value = "September 2022 (through September 15)"

# --- Update/insert the daily counts affected by the new day of data ---

# Austin (row 2)
$ws.Range("K2").Value = 6
$ws.Range("T2").Value = 4
$ws.Range("AL2").Value = 5
$ws.Range("AU2").Value = 3

# Englewood (row 6)
$ws.Range("D6").Value = 8

# Little Italy, UIC (row 7)
$ws.Range("BD7").Value = 1

# Humboldt Park (row 8)
$ws.Range("K8").Value = 3

# North Lawndale (row 10)
$ws.Range("T10").Value = 7

# Grand Crossing (row 12)
$ws.Range("B12").Value = 2
$ws.Range("D12").Value = 9
$ws.Range("T12").Value = 4

# Roseland (row 14)
$ws.Range("AL14").Value = 2

# Chatham (row 15)
$ws.Range("K15").Value = 5

# South Chicago (row 23)
$ws.Range("B23").Value = 2

# South Shore (row 24)
$ws.Range("AC24").Value = 1

# Calumet Heights (row 25)
$ws.Range("AC25").Value = 2

# Avondale (row 43)
$ws.Range("K43").Value = 1
$ws.Range("BD43").Value = 1

# Armour Square (row 60)
$ws.Range("T60").Value = 1

# Hermosa (row 73)
$ws.Range("K73").Value = 1

# Old Town (row 90)
$ws.Range("K90").Value = 1

# South Deering (row 95)
$ws.Range("AU95").Value = 1
